$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title paragraph: "Section 1: Week 1: Evaluate Cybersecurity"
#    -> "Section 1: Week 3: Global Security Risk"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Evaluate Cybersecurity", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Global Security Risk", 2)
$d.Content.Find.Execute("Week 1", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Week 3", 2)

# ------------------------------------------------------------------
# 2. Date paragraph: "April 5, 2020" -> "April 19, 2020"
# ------------------------------------------------------------------
$d.Content.Find.Execute("April 5, 2020", $true, $false, $false, $false, $false, `
    $true, 1, $false, "April 19, 2020", 2)

# ------------------------------------------------------------------
# 3. Body heading: "Evaluate CyberSecurity" -> "Global Security Risks"
#    and drop the leftover _GoBack bookmark.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Evaluate CyberSecurity", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Global Security Risks", 2)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 4. Append the paper outline as a multilevel numbered list.
# ------------------------------------------------------------------
$items = @(
    @{lvl=0; text="Choose the organization"},
    @{lvl=1; text="Ride-Me the ride-hailing app, billion-dollar valuation"},
    @{lvl=1; text="Domestically strong and seeking international growth"},
    @{lvl=1; text="Competitors like Uber, Lyft, and Google"},
    @{lvl=0; text="Section I: Understanding Risk"},
    @{lvl=1; text="Risks from the International community"},
    @{lvl=2; text="Espionage"},
    @{lvl=2; text="Sabotage"},
    @{lvl=2; text="Subversion"},
    @{lvl=1; text="Process for establishing risk"},
    @{lvl=2; text="Geography"},
    @{lvl=2; text="Legal Challenges"},
    @{lvl=2; text="National Sovereignty"},
    @{lvl=1; text="Threats"},
    @{lvl=2; text="State-sponsored actors"},
    @{lvl=2; text="Zero-day attacks"},
    @{lvl=0; text="Section II: Mitigating Risk"},
    @{lvl=1; text="People"},
    @{lvl=1; text="Process"},
    @{lvl=1; text="Products"},
    @{lvl=0; text="Section III: Budgeting Resources"},
    @{lvl=1; text="What does it take/cost"},
    @{lvl=1; text="Communicating Necessity"},
    @{lvl=0; text="Conclusion"}
)

$listTemplate = $null
$isFirstItem = $true

foreach ($item in $items) {
    $lastPara = $d.Paragraphs($d.Paragraphs.Count)
    $newPara = $d.Paragraphs.Add($lastPara.Range)
    $newPara.Range.Text = $item.text
    $newPara.Style = "ListParagraph"

    if ($isFirstItem) {
        # Seed the multilevel numbered list (numId=1) and shape it into
        # the classic 1./a./i. outline pattern (repeats every 3 levels).
        $newPara.Range.ListFormat.ApplyNumberDefault()
        $listTemplate = $newPara.Range.ListFormat.ListTemplate

        $numberStyles = @(0, 4, 2, 0, 4, 2, 0, 4, 2)
        for ($lvlIdx = 1; $lvlIdx -le 9; $lvlIdx++) {
            $listLevel = $listTemplate.ListLevels($lvlIdx)
            $listLevel.NumberStyle = $numberStyles[$lvlIdx - 1]
        }

        $isFirstItem = $false
    }
    else {
        $newPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
    }

    $newPara.Range.ListFormat.ListLevelNumber = $item.lvl + 1
}
